# Update the repayment strategy value on the ProductLoanInput sheet
# (row 17, column B) from "RBI (India)" to "Overdue/Due Fee/Int,Principal".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Reflect the selected cell used while making the edit.
$ws.Activate()
$ws.Range("B17").Select()
